$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.73
$ws.Range("I2").Value = 4.1
$ws.Range("J2").Value = 2.25
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 4
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("AJ2").Value = 13
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 26
$ws.Range("AT2").Value = 3.75
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 3.6
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("Z4").Value = 26
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 7
$ws.Range("AN4").Value = 4.5
$ws.Range("AT4").Value = 2.38
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.2
$ws.Range("L5").Value = 2.88
$ws.Range("N5").Value = 9.5
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 2
$ws.Range("X5").Value = 17
$ws.Range("AG5").Value = 201
$ws.Range("AJ5").Value = 9
$ws.Range("AO5").Value = 19
$ws.Range("AY5").Value = 12
$ws.Range("AZ5").Value = 21
$ws.Range("BB5").Value = 51
$ws.Range("G6").Value = 1.8
$ws.Range("I6").Value = 4.33
$ws.Range("J6").Value = 2.4
$ws.Range("L6").Value = 4.5
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 9
$ws.Range("Z6").Value = 15
$ws.Range("AE6").Value = 15
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 15
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 41
$ws.Range("AN6").Value = 3.75
$ws.Range("AO6").Value = 9.5
$ws.Range("AU6").Value = 8
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 29
$ws.Range("BA6").Value = 81
$ws.Range("BC6").Value = 201
